{"js": "// Update the benchmark stats table (single-column table, 1 cell per row).\n// Several rows get their numeric value replaced, and three rows that used\n// to hold a full tab-separated stats line are collapsed down to the single\n// summary value that used to live a few rows above them.\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// 0-based (row index, column index) -> new cell text\nconst updates = [\n  [0, 0, \"0M\"],\n  [1, 0, \"0M\"],\n  [2, 0, \"0M\"],\n  [3, 0, \"6476\"],\n  [5, 0, \"0.34009\"],\n  [6, 0, \"0.09437\"],\n  [7, 0, \"0.00509\"],\n  [8, 0, \"0.34009\"],\n  [9, 0, \"0.34009\"],\n  [10, 0, \"0.34009\"],\n  [11, 0, \"16.75468\"],\n  [43, 0, \"99.05\"],\n  [44, 0, \"16.75\"],\n  [45, 0, \"1764\"],\n];\n\nfor (const [row, col, text] of updates) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the benchmark stats table (single-column table, 1 cell per row).\n# Several rows get their numeric value replaced, and three rows that used\n# to hold a full tab-separated stats line are collapsed down to the single\n# summary value that used to live a few rows above them.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$updates = @{\n    1  = \"0M\"\n    2  = \"0M\"\n    3  = \"0M\"\n    4  = \"6476\"\n    6  = \"0.34009\"\n    7  = \"0.09437\"\n    8  = \"0.00509\"\n    9  = \"0.34009\"\n    10 = \"0.34009\"\n    11 = \"0.34009\"\n    12 = \"16.75468\"\n    44 = \"99.05\"\n    45 = \"16.75\"\n    46 = \"1764\"\n}\n\nforeach ($rowIndex in $updates.Keys) {\n    $cell = $t.Cell($rowIndex, 1)\n    $cell.Range.Text = $updates[$rowIndex]\n}\n"}
